$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the "Template" sheet to create "Week3", placed immediately
#    before "Template" (Week1, Week2, Week3, Template).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("Template")
$template.Copy($template)
$week3 = $wb.Worksheets.Item("Template (2)")
$week3.Name = "Week3"

# ---------------------------------------------------------------------------
# 2. Week2 ("Friday" column, column F) - fill in the diary entries.
#    Order matters: it determines the shared-string table order.
# ---------------------------------------------------------------------------
$week2 = $wb.Worksheets.Item("Week2")

$week2.Range("F5").Value = "Meeting with Joh Monday 3.30" + [char]10 + "Detailed flow of information"
$week2.Range("F3").Value = "Updated how the Phase1 may need to be structured. Have problem now with sending lots of requests. Meeting with john to confirm right course of action."
$week2.Range("F4").Value = "Asked people what they wanted out of QA in this organisation" + [char]10 + "Read Mcconnell chapter on quality assurance" + [char]10 + "Went through Microsoft tutorials and said which ones I could do (content)"
$week2.Range("F7").Value = "Getting over a new problem and finding new solutions to it"

# F4 got its wrap-text formatting switched on when the multi-line text was
# typed in (matches the other filled-in cells in that row).
$week2.Range("F4").WrapText = $true

# New column F got resized to fit the new content.
$week2.Columns.Item(6).ColumnWidth = 19.5703125

# ---------------------------------------------------------------------------
# 3. Week3 - fill in the diary entries that were typed into the new sheet.
# ---------------------------------------------------------------------------
$week3.Range("B4").Value = "Read up on book. Started looking for more books"
$week3.Range("B6").Value = "Create propper burn down chart + start propper planning after Johns advice"
$week3.Range("B7").Value = "Confidence in going into Johns office knowing what I am talking about and understanding what he is going on about"

$week3.Rows.Item(4).RowHeight = 37.5
$week3.Rows.Item(6).RowHeight = 63
$week3.Rows.Item(7).RowHeight = 75.75

$week3.Columns.Item(1).ColumnWidth = 24.5703125
$week3.Columns.Item(2).ColumnWidth = 23.28515625

# ---------------------------------------------------------------------------
# 4. View state - selections on each sheet, and which sheet/cell is active.
# ---------------------------------------------------------------------------
$week1 = $wb.Worksheets.Item("Week1")
$week1.Activate()
$week1.Range("B8").Select()

$week2.Activate()
$week2.Range("F6").Select()

$tmpl = $wb.Worksheets.Item("Template")
$tmpl.Activate()
$tmpl.Cells.Select()

# Week3 is the sheet left selected/active when the file was saved.
$week3.Activate()
$week3.Range("E16").Select()
